# mandate_tracker.xlsx -- "added comments to did code"
#
# This script reproduces, via Excel COM interop, the set of edits captured
# in the target OOXML diff:
#   * several "Year Expanded..." (column D) numeric-year cells are replaced
#     with the literal "NA" (treatment-year info moved to a Notes/Source
#     column instead)
#   * a handful of Source-column (G) cells that already held a URL as text
#     get turned into real clickable hyperlinks (style upgraded to the
#     existing wrap-text "Hyperlink" cell style)
#   * a couple of brand-new Source/Notes cells are added (Georgia source
#     link, "2017 potential treatment year" note)
#   * the West Virginia note that was split across two cells (F50/F51) is
#     merged back into a single F50 cell
#   * a couple of row heights grow to fit the (now) longer wrapped text
#   * the active selection moves to B20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: turn an existing text-only URL cell into a real hyperlink while
# preserving the workbook's existing wrap-text "Hyperlink" style (rather
# than letting Excel mint a brand-new style entry for it). We do this by
# adding the hyperlink first, then re-pasting the cell formats from a cell
# that already uses that exact style (G17), which snaps the cell back onto
# the pre-existing style slot.
# ---------------------------------------------------------------------
function Add-SourceHyperlink {
    param(
        [string]$CellRef,
        [string]$Url
    )
    $cell = $ws.Range($CellRef)
    if ($Url -match '^(.*)#(.*)$') {
        $base = $Matches[1]
        $frag = $Matches[2]
        $ws.Hyperlinks.Add($cell, $base, $frag) | Out-Null
    } else {
        $ws.Hyperlinks.Add($cell, $Url) | Out-Null
    }
    $ws.Range("G17").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Row 11 (Florida) -- taller row + G11 becomes a hyperlink
# ---------------------------------------------------------------------
$ws.Rows(11).RowHeight = 30
Add-SourceHyperlink "G11" $ws.Range("G11").Value2

# ---------------------------------------------------------------------
# Row 12 (Georgia) -- new source link
# ---------------------------------------------------------------------
$ws.Range("G12").Value2 = "https://www.medicaid.gov/Medicaid/downloads/georgia-mcp.pdf"

# ---------------------------------------------------------------------
# Row 15 (Illinois) -- taller row + G15 becomes a hyperlink
# ---------------------------------------------------------------------
$ws.Rows(15).RowHeight = 45
Add-SourceHyperlink "G15" $ws.Range("G15").Value2

# ---------------------------------------------------------------------
# Row 19 (Kentucky) -- G19 becomes a hyperlink
# ---------------------------------------------------------------------
Add-SourceHyperlink "G19" $ws.Range("G19").Value2

# ---------------------------------------------------------------------
# Row 26 (Mississippi) -- G26 becomes a hyperlink
# ---------------------------------------------------------------------
Add-SourceHyperlink "G26" $ws.Range("G26").Value2

# ---------------------------------------------------------------------
# Row 29 (Nebraska) -- D29 numeric year -> "NA"; new F29 note
# ---------------------------------------------------------------------
$ws.Range("D29").Value2 = "NA"
$ws.Range("F29").Value2 = "2017 potential treatment year"

# ---------------------------------------------------------------------
# Row 30 (Nevada) -- new D30 "NA"
# ---------------------------------------------------------------------
$ws.Range("D30").Value2 = "NA"

# ---------------------------------------------------------------------
# Row 34 (New York) -- new D34 "NA"
# ---------------------------------------------------------------------
$ws.Range("D34").Value2 = "NA"

# ---------------------------------------------------------------------
# Row 35 (North Carolina) -- G35 becomes a hyperlink (has a #:~:text= frag)
# ---------------------------------------------------------------------
Add-SourceHyperlink "G35" $ws.Range("G35").Value2

# ---------------------------------------------------------------------
# Row 37 (Ohio) -- G37 becomes a hyperlink
# ---------------------------------------------------------------------
Add-SourceHyperlink "G37" $ws.Range("G37").Value2

# ---------------------------------------------------------------------
# Row 42 (South Carolina) -- D42 numeric year -> "NA"; G42 becomes a link
# ---------------------------------------------------------------------
$ws.Range("D42").Value2 = "NA"
Add-SourceHyperlink "G42" $ws.Range("G42").Value2

# ---------------------------------------------------------------------
# Row 45 (Texas) -- G45 becomes a hyperlink
# ---------------------------------------------------------------------
Add-SourceHyperlink "G45" $ws.Range("G45").Value2

# ---------------------------------------------------------------------
# Row 46 (Utah) -- D46 numeric year -> "NA"; G46 becomes a link (has frag)
# ---------------------------------------------------------------------
$ws.Range("D46").Value2 = "NA"
Add-SourceHyperlink "G46" $ws.Range("G46").Value2

# ---------------------------------------------------------------------
# Row 48 (Virginia) -- new D48 "NA"
# ---------------------------------------------------------------------
$ws.Range("D48").Value2 = "NA"

# ---------------------------------------------------------------------
# Row 49 (Washington) -- G49 becomes a hyperlink
# ---------------------------------------------------------------------
Add-SourceHyperlink "G49" $ws.Range("G49").Value2

# ---------------------------------------------------------------------
# Row 50/51 (West Virginia / Wisconsin) -- merge the WV note that had
# spilled into F51 ("programs.") back into F50, taller row 50, and both
# D50/D51 numeric years -> "NA"
# ---------------------------------------------------------------------
$ws.Rows(50).RowHeight = 60
$f50 = $ws.Range("F50").Value2
$f51 = $ws.Range("F51").Value2
$ws.Range("F50").Value2 = $f50 + " " + $f51.TrimEnd(".")
$ws.Range("F51").ClearContents()

$ws.Range("D50").Value2 = "NA"
$ws.Range("D51").Value2 = "NA"

# ---------------------------------------------------------------------
# View state: move the active selection
# ---------------------------------------------------------------------
$ws.Range("B20").Select() | Out-Null
